$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.Range("I2:I135")
$used.Replace("phi.u", "phiu")
$used.Replace("eta.fu", "etafu")
$ws.Range("M22").Select()
